# Insert a new data row at row 473 on Sheet1, pushing all the existing
# rows (473..572) down by one (to 474..573), and populate the newly
# inserted row with the new "Ajo" (garlic) price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 473; this shifts every
# row from 473 downward one position (old 473 -> new 474, ..., old 572 -> new 573)
# and extends the sheet dimension from A1:R572 to A1:R573 automatically.
$ws.Rows.Item(473).Insert()

# Populate the newly inserted row 473 with the new record's values.
$ws.Range("A473").Value = 5
$ws.Range("B473").Value = "Macroferia Regional de Talca"
$ws.Range("C473").Value = "Maule"
$ws.Range("D473").Value = 45258
$ws.Range("E473").Value = 7
$ws.Range("F473").Value = 100112003
$ws.Range("G473").Value = "Ajo"
$ws.Range("H473").Value = "Chino"
$ws.Range("I473").Value = "Primera"
$ws.Range("J473").Value = 200
$ws.Range("K473").Value = 23000
$ws.Range("L473").Value = 23000
$ws.Range("M473").Value = 23000
$ws.Range("N473").Value = "$/malla 10 kilos"
$ws.Range("O473").Value = "Región del Maule"
$ws.Range("P473").Value = 2300
$ws.Range("Q473").Value = 10
$ws.Range("R473").Value = "Hortaliza"
